# Daily refresh of the cryptos worksheet (prices / 1h volume %) as produced
# by the scheduled GitHub Actions job. Column D ("Price") holds numeric-
# looking text (e.g. "0.999", "604.68") that must stay stored as text, just
# like in the source file, instead of being auto-coerced to a number by
# Excel's usual Range.Value type inference. Set-TextValue forces the cell to
# Text format for the assignment and then clears the format again so the
# cell ends up with no explicit style, matching the original cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Cells.Item(2, 4).Value = '68.958.73'
$ws.Cells.Item(2, 5).Value = '  -0.29%  '
$ws.Cells.Item(3, 4).Value = '3.936.74'
$ws.Cells.Item(3, 5).Value = '  +3.38%  '
Set-TextValue $ws.Cells.Item(4, 4) '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
Set-TextValue $ws.Cells.Item(5, 4) '604.68'
$ws.Cells.Item(5, 5).Value = '  +0.62%  '
Set-TextValue $ws.Cells.Item(6, 4) '167.90'
$ws.Cells.Item(6, 5).Value = '  +1.89%  '
$ws.Cells.Item(7, 4).Value = '3.937.09'
$ws.Cells.Item(7, 5).Value = '  +3.50%  '
$ws.Cells.Item(8, 5).Value = '  +0.10%  '
$ws.Cells.Item(9, 5).Value = '  -0.06%  '
$ws.Cells.Item(10, 5).Value = '  +0.84%  '
$ws.Cells.Item(11, 5).Value = '  +3.04%  '
$ws.Cells.Item(12, 5).Value = '  +1.56%  '
Set-TextValue $ws.Cells.Item(13, 4) '0.0000255'
$ws.Cells.Item(13, 5).Value = '  +3.88%  '
Set-TextValue $ws.Cells.Item(14, 4) '37.67'
$ws.Cells.Item(14, 5).Value = '  +1.47%  '
$ws.Cells.Item(15, 4).Value = '4.594.91'
$ws.Cells.Item(15, 5).Value = '  +3.39%  '
$ws.Cells.Item(16, 4).Value = '3.950.42'
$ws.Cells.Item(16, 5).Value = '  +4.03%  '
$ws.Cells.Item(17, 4).Value = '69.016.18'
$ws.Cells.Item(17, 5).Value = '  -0.37%  '
$ws.Cells.Item(18, 5).Value = '  -0.33%  '
Set-TextValue $ws.Cells.Item(19, 4) '17.44'
$ws.Cells.Item(19, 5).Value = '  +1.17%  '
$ws.Cells.Item(20, 5).Value = '  -1.81%  '
Set-TextValue $ws.Cells.Item(21, 4) '10.98'
$ws.Cells.Item(21, 5).Value = '  -3.47%  '
Set-TextValue $ws.Cells.Item(22, 4) '494.43'
$ws.Cells.Item(22, 5).Value = '  +1.39%  '
Set-TextValue $ws.Cells.Item(23, 4) '0.732'
$ws.Cells.Item(23, 5).Value = '  +1.69%  '
Set-TextValue $ws.Cells.Item(24, 4) '0.0000168'
$ws.Cells.Item(24, 5).Value = '  +6.64%  '
Set-TextValue $ws.Cells.Item(25, 4) '85.02'
$ws.Cells.Item(26, 5).Value = '  +1.02%  '
$ws.Cells.Item(27, 5).Value = '  +0.32%  '
$ws.Cells.Item(28, 5).Value = '  +1.88%  '
$ws.Cells.Item(29, 5).Value = '  +0.02%  '
Set-TextValue $ws.Cells.Item(30, 4) '2.99'
$ws.Cells.Item(30, 5).Value = '  +0.69%  '
$ws.Cells.Item(31, 4).Value = '4.088.52'
$ws.Cells.Item(31, 5).Value = '  +3.23%  '
$ws.Cells.Item(32, 5).Value = '  +0.13%  '
Set-TextValue $ws.Cells.Item(33, 4) '7.81'
$ws.Cells.Item(33, 5).Value = '  -2.53%  '
Set-TextValue $ws.Cells.Item(34, 4) '32.10'
$ws.Cells.Item(34, 5).Value = '  +1.17%  '
$ws.Cells.Item(35, 4).Value = '3.899.98'
$ws.Cells.Item(35, 5).Value = '  +3.97%  '
Set-TextValue $ws.Cells.Item(36, 4) '0.107'
$ws.Cells.Item(36, 5).Value = '  +0.31%  '
$ws.Cells.Item(37, 5).Value = '  +1.51%  '
$ws.Cells.Item(38, 5).Value = '  +2.66%  '
$ws.Cells.Item(39, 2).Value = 'dogwifhat'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Cells.Item(39, 4) '3.34'
$ws.Cells.Item(39, 5).Value = '  +9.87%  '
$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Cells.Item(40, 4) '0.139'
$ws.Cells.Item(40, 5).Value = '  -0.13%  '
$ws.Cells.Item(41, 5).Value = '  +0.06%  '
Set-TextValue $ws.Cells.Item(42, 4) '0.322'
$ws.Cells.Item(42, 5).Value = '  +0.88%  '
Set-TextValue $ws.Cells.Item(43, 4) '438.58'
$ws.Cells.Item(43, 5).Value = '  +0.45%  '
$ws.Cells.Item(44, 5).Value = '  +0.89%  '
Set-TextValue $ws.Cells.Item(45, 4) '48.16'
$ws.Cells.Item(45, 5).Value = '  -0.81%  '
Set-TextValue $ws.Cells.Item(46, 4) '8.62'
$ws.Cells.Item(46, 5).Value = '  +2.78%  '
Set-TextValue $ws.Cells.Item(48, 4) '0.000278'
$ws.Cells.Item(48, 5).Value = '  +22.74%  '
Set-TextValue $ws.Cells.Item(49, 4) '143.07'
$ws.Cells.Item(49, 5).Value = '  +0.49%  '
$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Cells.Item(50, 4) '0.0361'
$ws.Cells.Item(50, 5).Value = '  +2.09%  '
$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).Value = '2.820.21'
$ws.Cells.Item(51, 5).Value = '  -0.12%  '
